# "Implemented VS Studio Integration with Unity, Updated Story"
#
# Updates the Group-Planning task sheet:
#  - Row 14 ("Add Secret Sound for Cave Secret + Compose") is finished ->
#    mark it visually completed (strikethrough), same way rows 2-4/13 look.
#  - Rows 18 / 20 / 23 get their "Effort" (worked hours) filled in, which
#    finishes them off (Remain -> 0, Completion -> 100%).
#  - Row 28 gets extra curr. estimate + some effort logged (25% done).
#  - Rows 29-32 are brand-new backlog/finished items that get their story
#    text + numbers filled in.
#  - Selection cursor ends up on I32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: task is complete -> strike it through like the other
#     "done" rows (2-4, 13) which already use the struck-through style.
$ws.Range("A14:I14").Font.Strikethrough = $true

# --- Row 18: effort catches up to the current estimate (task finished)
$ws.Range("D18").Value = 9

# --- Row 20: effort logged equal to current estimate (task finished)
$ws.Range("D20").Value = 10

# --- Row 23: effort logged equal to current estimate (task finished)
$ws.Range("D23").Value = 3

# --- Row 28: current estimate revised up, some effort logged
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 1

# --- Row 29: new task - "Golem Sprite"
$ws.Range("A29").Value = "Golem Sprite"
$ws.Range("B29").Value = 2.5
$ws.Range("C29").Value = 2.5
$ws.Range("F29").Value = "Steffi"

# --- Row 30: new task - "End Sequence for first Level"
$ws.Range("A30").Value = "End Sequence for first Level"
$ws.Range("B30").Value = 4
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 3
$ws.Range("F30").Value = "Sascha"

# --- Row 31: new task - "Rotation of the player by Mouse Movement Fix"
$ws.Range("A31").Value = "Rotation of the player by Mouse Movement Fix"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 1
$ws.Range("F31").Value = "Cedric"

# --- Row 32: new task - "SquirelWurf Skript"
$ws.Range("A32").Value = "SquirelWurf Skript"
$ws.Range("B32").Value = 3
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 2
$ws.Range("F32").Value = "Cedric"

# --- Move the selection cursor to I32 (where the author's edits ended)
$null = $ws.Range("I32").Select()
